$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column R (18), shifting "Repeat HV" to column S.
$ws.Columns.Item(18).Insert()

# Set the header for the newly inserted column R2.
$ws.Range("R2").Value = "Unidad de medida"

# Match the column width Excel applies when inserting next to column Q.
$ws.Columns.Item(18).ColumnWidth = $ws.Columns.Item(17).ColumnWidth

# Update the selection/active cell to mirror the saved workbook state.
$ws.Range("R2").Select()
